$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "notas" value in A2
$ws.Range("A2").Value = 2000056514

# Move the active selection to A3 (was C7)
$ws.Range("A3").Select()

# Add the new (empty) result-of-effort cells J3:J6, formatted like column A
# (style index 2: centered alignment) so the used range grows to A1:J6
$ws.Range("J3:J6").HorizontalAlignment = -4108
